# Update the "Fixture" worksheet:
#  - Remove the first two rows (Wycombe Wanderers Carabao Cup match, and the
#    PSG Champions League away match) which have already been played.
#  - Insert a new fixture row for the West Ham United away match after the
#    Burnley fixture.
#  - Update the kick-off date/time for the home fixtures against Everton and
#    West Ham United, which were rescheduled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two fixtures that have already been played (old rows 1 and 2).
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# After the deletion, "Manchester City v Burnley" is now row 1. Insert a new
# row after it for the West Ham United away fixture.
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "West Ham United v Manchester City"
$ws.Range("B2").Value = "27/10/2021 19:45 | "

# Update rescheduled kick-off times for the Everton (row 4) and West Ham
# United (row 5) home fixtures.
$ws.Range("B4").Value = "21/11/2021 14:00 | Premier League"
$ws.Range("B5").Value = "28/11/2021 14:00 | Premier League"
